# "refine hotel voting use case"
# Moves the vote checkmark ("√") within two of the hotel-voting blocks:
#   - Day 2 block: checkmark moves from A18 (Globetrotter Hotel) to A20 (LJK Warrens Resort)
#   - Day 3 block: checkmark moves from A26 (Grandiose Estate Hotel) to A27 (Viewpoint Resort)
# and updates the sheet selection to A14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("voting")

# Day 2 block: clear old vote mark, set new one
$ws.Range("A18").Value = ""
$ws.Range("A20").Value = "√"

# Day 3 block: clear old vote mark, set new one
$ws.Range("A26").Value = ""
$ws.Range("A27").Value = "√"

# Update the active selection shown when the sheet is reopened
[void]$ws.Range("A14").Select()
